# Apply cryptocurrency price/volume updates for Wed Aug  2 06:17:41 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextCell "D2" "29.645.14"
$ws.Range("E2").Value = "  +2.51%  "

Set-TextCell "D3" "1.861.41"
$ws.Range("E3").Value = "  +1.75%  "

Set-TextCell "D4" "0.9991"
$ws.Range("E4").Value = "  -0.02%  "

Set-TextCell "D5" "245.07"
$ws.Range("E5").Value = "  +1.90%  "

Set-TextCell "D6" "0.6978"
$ws.Range("E6").Value = "  +1.49%  "

$ws.Range("E7").Value = "  +0.01%  "

Set-TextCell "D8" "0.07714"
$ws.Range("E8").Value = "  +0.96%  "

Set-TextCell "D9" "0.3062"
$ws.Range("E9").Value = "  +0.62%  "

Set-TextCell "D10" "23.71"
$ws.Range("E10").Value = "  +0.71%  "

Set-TextCell "D11" "0.07753"
$ws.Range("E11").Value = "  -0.35%  "

Set-TextCell "D12" "5.165"
$ws.Range("E12").Value = "  +2.03%  "

Set-TextCell "D13" "1.856.36"
$ws.Range("E13").Value = "  +1.57%  "

Set-TextCell "D14" "92.30"
$ws.Range("E14").Value = "  +2.05%  "

Set-TextCell "D15" "0.6931"
$ws.Range("E15").Value = "  +2.74%  "

Set-TextCell "D16" "6.572"
$ws.Range("E16").Value = "  +2.24%  "

Set-TextCell "D17" "29.631.77"
$ws.Range("E17").Value = "  +2.53%  "

Set-TextCell "D18" "0.000008322"
$ws.Range("E18").Value = "  +0.62%  "

Set-TextCell "D19" "2.104.33"
$ws.Range("E19").Value = "  +1.59%  "

Set-TextCell "D20" "241.32"
$ws.Range("E20").Value = "  -0.49%  "

$ws.Range("E21").Value = "  +0.88%  "

Set-TextCell "D22" "0.9994"
$ws.Range("E22").Value = "  -0.04%  "

Set-TextCell "D23" "7.618"
$ws.Range("E23").Value = "  +2.77%  "

Set-TextCell "D24" "0.9998"
$ws.Range("E24").Value = "  +0.08%  "

Set-TextCell "D25" "0.1502"
$ws.Range("E25").Value = "  +1.90%  "

Set-TextCell "D26" "8.928"
$ws.Range("E26").Value = "  +1.75%  "

Set-TextCell "D27" "159.19"
$ws.Range("E27").Value = "  -1.35%  "

$ws.Range("E28").Value = "  +0.58%  "

Set-TextCell "D29" "1.536"
$ws.Range("E29").Value = "  +0.04%  "

Set-TextCell "D30" "4.257"
$ws.Range("E30").Value = "  +1.15%  "

Set-TextCell "D31" "4.185"
$ws.Range("E31").Value = "  +1.37%  "

Set-TextCell "D32" "1.197"
$ws.Range("E32").Value = "  +0.55%  "

Set-TextCell "D33" "0.05092"
$ws.Range("E33").Value = "  -0.36%  "

Set-TextCell "D34" "0.7770"
$ws.Range("E34").Value = "  +3.99%  "

Set-TextCell "D35" "1.901"
$ws.Range("E35").Value = "  +4.44%  "

$ws.Range("E36").Value = "  +1.13%  "

Set-TextCell "D37" "2.685"
$ws.Range("E37").Value = "  +0.31%  "

Set-TextCell "D38" "1.324.79"
$ws.Range("E38").Value = "  +8.92%  "

$ws.Range("E39").Value = "  +1.81%  "

Set-TextCell "D40" "2.733"
$ws.Range("E40").Value = "  +2.01%  "

Set-TextCell "D41" "0.9734"
$ws.Range("E41").Value = "  +6.27%  "

Set-TextCell "D42" "106.88"
$ws.Range("E42").Value = "  -1.36%  "

Set-TextCell "D43" "5.838"
$ws.Range("E43").Value = "  +9.56%  "

Set-TextCell "D44" "0.9993"
$ws.Range("E44").Value = "  -0.01%  "

$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextCell "D45" "0.00000000126"
$ws.Range("E45").Value = "  +3.72%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell "D46" "9.777"
$ws.Range("E46").Value = "  +3.14%  "

$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextCell "D47" "2.004.95"
$ws.Range("E47").Value = "  +1.59%  "

Set-TextCell "D48" "0.5214"
$ws.Range("E48").Value = "  +0.90%  "

Set-TextCell "D49" "1.780"
$ws.Range("E49").Value = "  +2.98%  "

Set-TextCell "D50" "63.90"
$ws.Range("E50").Value = "  +1.47%  "

Set-TextCell "D51" "6.973"
$ws.Range("E51").Value = "  +1.23%  "
